$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 64
$ws.Range("H64").Value = 7533
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 7533
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 7533
$ws.Range("M64").Value = ""
$ws.Range("N64").Value = -8029
# Row 67
$ws.Range("H67").Value = 7533
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 7533
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 7533
$ws.Range("M67").Value = ""
$ws.Range("N67").Value = -9249
# Row 86
$ws.Range("H86").Value = 5226.7896
$ws.Range("I86").Value = 4565.9
$ws.Range("J86").Value = 5961.1113
$ws.Range("K86").Value = 4565.9
$ws.Range("L86").Value = 5961.1113
$ws.Range("M86").Value = -3442.9
$ws.Range("N86").Value = -8207.1113
# Row 89
$ws.Range("H89").Value = 5226.7896
$ws.Range("I89").Value = 4565.9
$ws.Range("J89").Value = 5961.1113
$ws.Range("K89").Value = 22829.5
$ws.Range("L89").Value = 29805.5565
$ws.Range("M89").Value = -17213.5
$ws.Range("N89").Value = -41037.5565
# Row 118
$ws.Range("H118").Value = 528.7143
$ws.Range("I118").Value = 571
$ws.Range("J118").Value = 275
$ws.Range("K118").Value = 1713
$ws.Range("L118").Value = 825
$ws.Range("M118").Value = -56
$ws.Range("N118").Value = -4139
# Row 132
$ws.Range("H132").Value = 33340276
$ws.Range("I132").Value = 50007584
$ws.Range("J132").Value = 5660.9
$ws.Range("K132").Value = 150022752
$ws.Range("L132").Value = 16982.7
$ws.Range("M132").Value = -150020222
$ws.Range("N132").Value = -22042.7
# Row 137
$ws.Range("H137").Value = 2975.4707
$ws.Range("I137").Value = 2256.1428
$ws.Range("J137").Value = 3479
$ws.Range("K137").Value = 6768.428400000001
$ws.Range("L137").Value = 10437
$ws.Range("M137").Value = -4218.428400000001
$ws.Range("N137").Value = -15537
# Row 138
$ws.Range("H138").Value = 627983.0600000001
$ws.Range("I138").Value = 1830.3334
$ws.Range("J138").Value = 1003674.7
$ws.Range("K138").Value = 5491.0002
$ws.Range("L138").Value = 3011024.1
$ws.Range("M138").Value = -351.0002000000004
$ws.Range("N138").Value = -3021304.1

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 36
$ws.Range("H36").Value = 10388.833
$ws.Range("I36").Value = 7833.25
$ws.Range("J36").Value = 15500
$ws.Range("K36").Value = 7833.25
$ws.Range("L36").Value = 15500
$ws.Range("M36").Value = -7487.25
$ws.Range("N36").Value = -16192
# Row 45
$ws.Range("H45").Value = 5065.091
$ws.Range("I45").Value = 3964.5625
$ws.Range("J45").Value = 7999.8335
$ws.Range("K45").Value = 3964.5625
$ws.Range("L45").Value = 7999.8335
$ws.Range("M45").Value = -3587.5625
$ws.Range("N45").Value = -8753.833500000001
# Row 69
$ws.Range("H69").Value = 329999
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 329999
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 329999
$ws.Range("M69").Value = ""
$ws.Range("N69").Value = -331497
# Row 72
$ws.Range("H72").Value = 329999
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 329999
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 989997
$ws.Range("M72").Value = ""
$ws.Range("N72").Value = -997485
# Row 122
$ws.Range("H122").Value = 4591.8423
$ws.Range("I122").Value = 4231.8
$ws.Range("J122").Value = 5439
$ws.Range("K122").Value = 12695.4
$ws.Range("L122").Value = 16317
$ws.Range("M122").Value = -10245.4
$ws.Range("N122").Value = -21217

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 42
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = ""
$ws.Range("N42").Value = ""
# Row 107
$ws.Range("H107").Value = 10420.571
$ws.Range("I107").Value = 12428
$ws.Range("J107").Value = 8413.143
$ws.Range("K107").Value = 12428
$ws.Range("L107").Value = 8413.143
$ws.Range("M107").Value = -10508
$ws.Range("N107").Value = -12253.143
# Row 130
$ws.Range("H130").Value = 64000
$ws.Range("I130").Value = 64000
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 64000
$ws.Range("L130").Value = 0
$ws.Range("M130").Value = -58980
$ws.Range("N130").Value = ""
# Row 134
$ws.Range("H134").Value = 2327.2354
$ws.Range("I134").Value = 1875.069
$ws.Range("J134").Value = 4949.8
$ws.Range("K134").Value = 5625.207
$ws.Range("L134").Value = 14849.4
$ws.Range("M134").Value = -3090.207
$ws.Range("N134").Value = -19919.4

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3089.3914
$ws.Range("I31").Value = 2023.25
$ws.Range("J31").Value = 4252.4546
$ws.Range("K31").Value = 2023.25
$ws.Range("L31").Value = 4252.4546
$ws.Range("M31").Value = -1728.25
$ws.Range("N31").Value = -4842.4546
# Row 34
$ws.Range("H34").Value = 3089.3914
$ws.Range("I34").Value = 2023.25
$ws.Range("J34").Value = 4252.4546
$ws.Range("K34").Value = 2023.25
$ws.Range("L34").Value = 4252.4546
$ws.Range("M34").Value = -1821.25
$ws.Range("N34").Value = -4656.4546
# Row 44
$ws.Range("H44").Value = 9666.333000000001
$ws.Range("I44").Value = 9666.333000000001
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 9666.333000000001
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = -9224.333000000001
$ws.Range("N44").Value = ""
# Row 58
$ws.Range("H58").Value = 2314.2354
$ws.Range("I58").Value = 1529.5
$ws.Range("J58").Value = 3435.2856
$ws.Range("K58").Value = 1529.5
$ws.Range("L58").Value = 3435.2856
$ws.Range("M58").Value = -1326.5
$ws.Range("N58").Value = -3841.2856
# Row 64
$ws.Range("H64").Value = 100000
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 100000
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 100000
$ws.Range("M64").Value = ""
$ws.Range("N64").Value = -100496
# Row 67
$ws.Range("H67").Value = 100000
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 100000
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 100000
$ws.Range("M67").Value = ""
$ws.Range("N67").Value = -101716
# Row 94
$ws.Range("H94").Value = 1876
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 1876
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 1876
$ws.Range("M94").Value = ""
$ws.Range("N94").Value = -2778
# Row 107
$ws.Range("H107").Value = 1064
$ws.Range("I107").Value = 810.3333
$ws.Range("J107").Value = 1419.1333
$ws.Range("K107").Value = 810.3333
$ws.Range("L107").Value = 1419.1333
$ws.Range("M107").Value = 1109.6667
$ws.Range("N107").Value = -5259.1333
# Row 129
$ws.Range("H129").Value = 25000
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 25000
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 25000
$ws.Range("M129").Value = ""
$ws.Range("N129").Value = -35000
# Row 132
$ws.Range("H132").Value = 857.8570999999999
$ws.Range("I132").Value = 834.1667
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 2502.5001
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = 27.4998999999998
$ws.Range("N132").Value = -8060
# Row 136
$ws.Range("H136").Value = 2314.2354
$ws.Range("I136").Value = 1529.5
$ws.Range("J136").Value = 3435.2856
$ws.Range("K136").Value = 4588.5
$ws.Range("L136").Value = 10305.8568
$ws.Range("M136").Value = -2038.5
$ws.Range("N136").Value = -15405.8568

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 140
$ws.Range("H140").Value = 3940.0952
$ws.Range("I140").Value = 2939.1538
$ws.Range("J140").Value = 5566.625
$ws.Range("K140").Value = 8817.4614
$ws.Range("L140").Value = 16699.875
$ws.Range("M140").Value = -3637.4614
$ws.Range("N140").Value = -27059.875

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 29
$ws.Range("H29").Value = 22832.666
$ws.Range("I29").Value = 20000
$ws.Range("J29").Value = 24249
$ws.Range("K29").Value = 20000
$ws.Range("L29").Value = 24249
$ws.Range("M29").Value = -19710
$ws.Range("N29").Value = -24829
# Row 80
$ws.Range("H80").Value = 6645.357
$ws.Range("I80").Value = 8483.286
$ws.Range("J80").Value = 4807.4287
$ws.Range("K80").Value = 8483.286
$ws.Range("L80").Value = 4807.4287
$ws.Range("M80").Value = -7485.286
$ws.Range("N80").Value = -6803.4287
# Row 83
$ws.Range("H83").Value = 6645.357
$ws.Range("I83").Value = 8483.286
$ws.Range("J83").Value = 4807.4287
$ws.Range("K83").Value = 42416.43
$ws.Range("L83").Value = 24037.1435
$ws.Range("M83").Value = -37424.43
$ws.Range("N83").Value = -34021.14350000001
# Row 101
$ws.Range("H101").Value = 54732.57
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 54732.57
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 54732.57
$ws.Range("M101").Value = ""
$ws.Range("N101").Value = -61222.57
# Row 128
$ws.Range("H128").Value = 87765
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 87765
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 87765
$ws.Range("M128").Value = ""
$ws.Range("N128").Value = -97725
# Row 140
$ws.Range("H140").Value = 87000
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 87000
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 87000
$ws.Range("M140").Value = ""
$ws.Range("N140").Value = -97360

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 82
$ws.Range("H82").Value = 14436.546
$ws.Range("I82").Value = 17983.666
$ws.Range("J82").Value = 10180
$ws.Range("K82").Value = 17983.666
$ws.Range("L82").Value = 10180
$ws.Range("M82").Value = -17622.666
$ws.Range("N82").Value = -10902
# Row 85
$ws.Range("H85").Value = 14436.546
$ws.Range("I85").Value = 17983.666
$ws.Range("J85").Value = 10180
$ws.Range("K85").Value = 17983.666
$ws.Range("L85").Value = 10180
$ws.Range("M85").Value = -16735.666
$ws.Range("N85").Value = -12676
# Row 132
$ws.Range("H132").Value = 5820.136
$ws.Range("I132").Value = 5126.6924
$ws.Range("J132").Value = 6821.778
$ws.Range("K132").Value = 15380.0772
$ws.Range("L132").Value = 20465.334
$ws.Range("M132").Value = -12850.0772
$ws.Range("N132").Value = -25525.334
# Row 134
$ws.Range("H134").Value = 126665
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 126665
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 126665
$ws.Range("M134").Value = ""
$ws.Range("N134").Value = -136805
# Row 136
$ws.Range("H136").Value = 3615.6316
$ws.Range("I136").Value = 3600.147
$ws.Range("J136").Value = 3747.25
$ws.Range("K136").Value = 10800.441
$ws.Range("L136").Value = 11241.75
$ws.Range("M136").Value = -8250.440999999999
$ws.Range("N136").Value = -16341.75

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 198583
$ws.Range("I62").Value = 262925.66
$ws.Range("J62").Value = 5555
$ws.Range("K62").Value = 262925.66
$ws.Range("L62").Value = 5555
$ws.Range("M62").Value = -262301.66
$ws.Range("N62").Value = -6803
# Row 65
$ws.Range("H65").Value = 198583
$ws.Range("I65").Value = 262925.66
$ws.Range("J65").Value = 5555
$ws.Range("K65").Value = 1314628.3
$ws.Range("L65").Value = 27775
$ws.Range("M65").Value = -1311508.3
$ws.Range("N65").Value = -34015
# Row 100
$ws.Range("H100").Value = 992.3
$ws.Range("I100").Value = 820.05884
$ws.Range("J100").Value = 1968.3334
$ws.Range("K100").Value = 1640.11768
$ws.Range("L100").Value = 3936.6668
$ws.Range("M100").Value = -1099.11768
$ws.Range("N100").Value = -5018.6668
# Row 104
$ws.Range("H104").Value = 27730
$ws.Range("I104").Value = 24420
$ws.Range("J104").Value = 28833.334
$ws.Range("K104").Value = 24420
$ws.Range("L104").Value = 28833.334
$ws.Range("M104").Value = -20926
$ws.Range("N104").Value = -35821.334
# Row 113
$ws.Range("H113").Value = 4167937.8
$ws.Range("I113").Value = 13889654
$ws.Range("J113").Value = 1487.9286
$ws.Range("K113").Value = 41668962
$ws.Range("L113").Value = 4463.7858
$ws.Range("M113").Value = -41666792
$ws.Range("N113").Value = -8803.7858
# Row 119
$ws.Range("H119").Value = 20349
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 20349
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 20349
$ws.Range("M119").Value = ""
$ws.Range("N119").Value = -30025
# Row 122
$ws.Range("H122").Value = 7015.3887
$ws.Range("I122").Value = 7311.6
$ws.Range("J122").Value = 6342.1816
$ws.Range("K122").Value = 21934.8
$ws.Range("L122").Value = 19026.5448
$ws.Range("M122").Value = -19484.8
$ws.Range("N122").Value = -23926.5448
# Row 132
$ws.Range("H132").Value = 1740.6
$ws.Range("I132").Value = 1234.6666
$ws.Range("J132").Value = 2499.5
$ws.Range("K132").Value = 3703.9998
$ws.Range("L132").Value = 7498.5
$ws.Range("M132").Value = -1173.9998
$ws.Range("N132").Value = -12558.5
